$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 249
$ws.Range("F7").Value = 12938
$ws.Range("F8").Value = 49
$ws.Range("F9").Value = 117
$ws.Range("F10").Value = 239
$ws.Range("F11").Value = 2998
$ws.Range("F13").Value = 6366
$ws.Range("F14").Value = 63
$ws.Range("F16").Value = 3368
$ws.Range("F18").Value = 160
$ws.Range("F19").Value = 121
$ws.Range("F20").Value = 36
$ws.Range("F21").Value = 57
$ws.Range("F23").Value = 27
$ws.Range("F24").Value = 3576
$ws.Range("F27").Value = 2735
$ws.Range("F28").Value = 2735
$ws.Range("F29").Value = 400
$ws.Range("F30").Value = 1870
$ws.Range("F31").Value = 100
$ws.Range("F32").Value = 207
$ws.Range("F33").Value = 6555
$ws.Range("F34").Value = 17
$ws.Range("F36").Value = 651
$ws.Range("F37").Value = 1967
$ws.Range("F38").Value = 1291
$ws.Range("F39").Value = 93
$ws.Range("F40").Value = 1026
$ws.Range("F41").Value = 15
$ws.Range("F43").Value = 218
$ws.Range("F46").Value = 125
$ws.Range("F47").Value = 1190
$ws.Range("F48").Value = 1750
$ws.Range("F49").Value = 153

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 32

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 423
$ws.Range("F3").Value = 584
$ws.Range("F4").Value = 13

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 32
$ws.Range("F6").Value = 423
$ws.Range("F7").Value = 584
$ws.Range("F8").Value = 249
$ws.Range("F10").Value = 12938
$ws.Range("F11").Value = 117
$ws.Range("F13").Value = 239
$ws.Range("F14").Value = 2998
$ws.Range("F15").Value = 6366
$ws.Range("F16").Value = 63
$ws.Range("F17").Value = 3368
$ws.Range("F19").Value = 160
$ws.Range("F20").Value = 121
$ws.Range("F21").Value = 36
$ws.Range("F22").Value = 57
$ws.Range("F25").Value = 27
$ws.Range("F26").Value = 3576
$ws.Range("F28").Value = 2735
$ws.Range("F29").Value = 400
$ws.Range("F30").Value = 1870
$ws.Range("F31").Value = 100
$ws.Range("F32").Value = 207
$ws.Range("F33").Value = 6555
$ws.Range("F35").Value = 17
$ws.Range("F37").Value = 652
$ws.Range("F38").Value = 1967
$ws.Range("F40").Value = 1291
$ws.Range("F41").Value = 93
$ws.Range("F42").Value = 1026
$ws.Range("F44").Value = 218
$ws.Range("F46").Value = 125
$ws.Range("F48").Value = 1751
$ws.Range("F50").Value = 153
